$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with new values (columns F and G for rows 306-357)
$ws.Cells.Item(306, 6).Value = 74350
$ws.Cells.Item(306, 7).Value = 7502
$ws.Cells.Item(307, 6).Value = 77900
$ws.Cells.Item(307, 7).Value = 6626
$ws.Cells.Item(309, 6).Value = 77411
$ws.Cells.Item(309, 7).Value = 5445
$ws.Cells.Item(311, 6).Value = 61657
$ws.Cells.Item(313, 6).Value = 73790
$ws.Cells.Item(317, 6).Value = 63733
$ws.Cells.Item(317, 7).Value = 2190
$ws.Cells.Item(322, 6).Value = 110385
$ws.Cells.Item(322, 7).Value = 2392
$ws.Cells.Item(323, 6).Value = 216162
$ws.Cells.Item(323, 7).Value = 3203
$ws.Cells.Item(324, 6).Value = 238089
$ws.Cells.Item(324, 7).Value = 2754
$ws.Cells.Item(327, 6).Value = 235254
$ws.Cells.Item(327, 7).Value = 2864
$ws.Cells.Item(328, 6).Value = 180139
$ws.Cells.Item(328, 7).Value = 2650
$ws.Cells.Item(329, 6).Value = 83333
$ws.Cells.Item(329, 7).Value = 1751
$ws.Cells.Item(330, 6).Value = 72434
$ws.Cells.Item(330, 7).Value = 2086
$ws.Cells.Item(331, 6).Value = 151341
$ws.Cells.Item(331, 7).Value = 2648
$ws.Cells.Item(332, 6).Value = 444716
$ws.Cells.Item(332, 7).Value = 4440
$ws.Cells.Item(333, 6).Value = 270713
$ws.Cells.Item(333, 7).Value = 2905
$ws.Cells.Item(334, 6).Value = 203112
$ws.Cells.Item(334, 7).Value = 3395
$ws.Cells.Item(335, 6).Value = 130251
$ws.Cells.Item(335, 7).Value = 2949
$ws.Cells.Item(336, 6).Value = 102335
$ws.Cells.Item(336, 7).Value = 3248
$ws.Cells.Item(337, 6).Value = 104516
$ws.Cells.Item(337, 7).Value = 2954
$ws.Cells.Item(338, 6).Value = 220850
$ws.Cells.Item(338, 7).Value = 3069
$ws.Cells.Item(339, 6).Value = 648718
$ws.Cells.Item(339, 7).Value = 5562
$ws.Cells.Item(340, 6).Value = 379847
$ws.Cells.Item(340, 7).Value = 3259
$ws.Cells.Item(341, 6).Value = 292303
$ws.Cells.Item(341, 7).Value = 3593
$ws.Cells.Item(342, 6).Value = 174663
$ws.Cells.Item(342, 7).Value = 2956
$ws.Cells.Item(343, 6).Value = 127409
$ws.Cells.Item(343, 7).Value = 2845
$ws.Cells.Item(344, 6).Value = 132021
$ws.Cells.Item(344, 7).Value = 2456
$ws.Cells.Item(345, 6).Value = 279441
$ws.Cells.Item(345, 7).Value = 3189
$ws.Cells.Item(346, 6).Value = 647820
$ws.Cells.Item(346, 7).Value = 4617
$ws.Cells.Item(347, 6).Value = 328948
$ws.Cells.Item(347, 7).Value = 2771
$ws.Cells.Item(348, 6).Value = 225301
$ws.Cells.Item(349, 6).Value = 152737
$ws.Cells.Item(349, 7).Value = 2625
$ws.Cells.Item(350, 6).Value = 122275
$ws.Cells.Item(350, 7).Value = 2663
$ws.Cells.Item(351, 6).Value = 141757
$ws.Cells.Item(351, 7).Value = 2637
$ws.Cells.Item(352, 6).Value = 289845
$ws.Cells.Item(352, 7).Value = 3404
$ws.Cells.Item(353, 6).Value = 669537
$ws.Cells.Item(353, 7).Value = 4911
$ws.Cells.Item(354, 6).Value = 284383
$ws.Cells.Item(354, 7).Value = 2614
$ws.Cells.Item(355, 6).Value = 208360
$ws.Cells.Item(355, 7).Value = 3192
$ws.Cells.Item(356, 6).Value = 150239
$ws.Cells.Item(356, 7).Value = 2699
$ws.Cells.Item(357, 6).Value = 125160
$ws.Cells.Item(357, 7).Value = 2747

# Add new row 358
$ws.Cells.Item(358, 1).Value = 44252
$ws.Cells.Item(358, 2).Value = 303420
$ws.Cells.Item(358, 3).Value = 12065
$ws.Cells.Item(358, 4).Value = 2645
$ws.Cells.Item(358, 5).Value = 6966
$ws.Cells.Item(358, 6).Value = 118387
$ws.Cells.Item(358, 7).Value = 2409
